$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; this shifts existing rows 53-77 down to 54-78
# (matches the target dimension change A1:R77 -> A1:R78)
$ws.Rows(53).Insert()

# Populate the newly inserted row 53 with the new record
$ws.Cells.Item(53, 1).Value = 3
$ws.Cells.Item(53, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(53, 3).Value = "Coquimbo"
$ws.Cells.Item(53, 4).Value = 44489
$ws.Cells.Item(53, 5).Value = 5
$ws.Cells.Item(53, 6).Value = 100112026
$ws.Cells.Item(53, 7).Value = "Haba"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 55
$ws.Cells.Item(53, 11).Value = 9000
$ws.Cells.Item(53, 12).Value = 9000
$ws.Cells.Item(53, 13).Value = 9000
$ws.Cells.Item(53, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(53, 16).Value = 360
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"
